# Add a new "Timeline" worksheet as the last sheet in the workbook and
# populate it with the weekly dissertation timeline.

$wb = $excel.ActiveWorkbook

# Insert the new sheet after the current last sheet so it lands at the end.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Timeline"

# Week number / date-range columns (C/D) plus task notes (column F).
$ws.Range("C4").Value = "Week 1"
$ws.Range("D4").Value = "17-23.6"

$ws.Range("C5").Value = "Week 2"
$ws.Range("D5").Value = "24-30.6"
$ws.Range("F5").Value = "Finish modularizing all the models"

$ws.Range("C6").Value = "Week 3"
$ws.Range("D6").Value = "1-7.7"
$ws.Range("F6").Value = "Think about the different kinds of data, gather them and organize"

$ws.Range("C7").Value = "Week 4"
$ws.Range("D7").Value = "8-14.7"
$ws.Range("F7").Value = "Run models on all data and record all results"

$ws.Range("C8").Value = "Week 5"
$ws.Range("D8").Value = "15-21.7"
$ws.Range("F8").Value = "Think about fine-tuning the models and hyperparameter selection"

$ws.Range("C9").Value = "Week 6"
$ws.Range("D9").Value = "22-28.7"
$ws.Range("F9").Value = "(potentially) create functionality to fine-tune the models"

$ws.Range("C10").Value = "Week 7"
$ws.Range("C11").Value = "Week 8"
$ws.Range("C12").Value = "Week 9"
$ws.Range("C13").Value = "Week 10"
$ws.Range("C14").Value = "Week 11"
$ws.Range("C15").Value = "Week 12"
$ws.Range("C16").Value = "Week 13"
$ws.Range("C17").Value = "Week 14"
$ws.Range("C18").Value = "Week 15"
$ws.Range("C19").Value = "Week 16"

# Match the author's final selection/view state on the new sheet.
[void]$ws.Range("F9").Select()
